# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) on each data sheet, de-bolds /
# un-borders the former "header-styled" row labels in column A (rows 2+),
# fixes a few accented-character typos, removes the now-unused "Teto" row
# on the Emissoes sheet, and refreshes the Custo Total sheet's header/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: strip the bold/border/center "header" look from a cell by
# copying plain number-only formats isn't needed here — these label
# cells only ever carry the Normal style once de-headered, so we just
# reset to the built-in "Normal" style (style index 0, i.e. no explicit
# style at all in the saved XML).
# ---------------------------------------------------------------------

function Clear-HeaderStyle($ws, $addr) {
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheets 1-4 share the exact same row layout/labels:
#   Potencia Acumulada - SIN (MW)
#   Geracao Periodo Medio (MWMed)
#   Atendimento a Ponta(MW)
#   Potencia Incremental - SIN(MW)
# Add an "Fonte/Tecnologia" header in A1 (copying the bold/border/center
# look already used by B1:E1 / A2), then drop the header look from the
# row-label cells A2:A12 and correct a few accents.
# ---------------------------------------------------------------------

$fonteSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $fonteSheets) {
    $ws = $wb.Worksheets.Item($name)

    # Grab the existing header look (bold, thin border, center/top align)
    # from A2 before we strip it, and paste it onto the new A1 header.
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    Clear-HeaderStyle $ws "A2"
    $ws.Range("A2").Value = "Hidro"

    Clear-HeaderStyle $ws "A3"
    $ws.Range("A3").Value = "Gás Natural"

    Clear-HeaderStyle $ws "A4"
    $ws.Range("A4").Value = "Carvão"

    Clear-HeaderStyle $ws "A5"
    $ws.Range("A5").Value = "Nuclear"

    Clear-HeaderStyle $ws "A6"
    $ws.Range("A6").Value = "Óleos Comb"

    Clear-HeaderStyle $ws "A7"
    $ws.Range("A7").Value = "Biomassa"

    Clear-HeaderStyle $ws "A8"
    $ws.Range("A8").Value = "Eólica"

    Clear-HeaderStyle $ws "A9"
    $ws.Range("A9").Value = "Solar"

    Clear-HeaderStyle $ws "A10"
    $ws.Range("A10").Value = "Outros"

    Clear-HeaderStyle $ws "A11"
    $ws.Range("A11").Value = "Pot. Compl."

    Clear-HeaderStyle $ws "A12"
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------
# Sheet 5: Emissoes Totais (MtCO2eq)
#   - add "Período" header in A1
#   - de-header / fix accents on A2, A3
#   - drop the now-unused row 4 ("Teto") entirely
# ---------------------------------------------------------------------

$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("A2").Copy() | Out-Null
$ws5.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws5.Range("A1").Value = "Período"

Clear-HeaderStyle $ws5 "A2"
$ws5.Range("A2").Value = "P.Médio"

Clear-HeaderStyle $ws5 "A3"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows.Item(4).Delete() | Out-Null

# ---------------------------------------------------------------------
# Sheet 6: Custo Total (bilhões de R$)
#   - add "Tipo Expansão" header in A1, B1 becomes "2015"
#   - de-header / fix accents on A2, A3 and update their values
# ---------------------------------------------------------------------

$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("B1").Copy() | Out-Null
$ws6.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws6.Range("A1").Value = "Tipo Expansão"

# B1's text needs to flip from "Custo" to the literal text "2015" (kept as
# a text string, not a number) — paste the already-textual "2015" value
# from another sheet's header row so it doesn't get auto-coerced to a number.
$wsFonte1 = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
$wsFonte1.Range("B1").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4163) | Out-Null   # xlPasteValues

Clear-HeaderStyle $ws6 "A2"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 690

Clear-HeaderStyle $ws6 "A3"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
